# Azure Batch and HPC.pptx edit
# Commit: "Put I|P|SaaS graphic on intro slides. Removed A4R slides from HPC & big data"
#
# The only real structural change is the removal of the blank "A4R" slide
# (original slide 2 - a title-only placeholder slide with no text at all)
# from the HPC deck. Every other slide keeps its own content/identity and
# simply shifts up one position as a natural consequence of the deletion:
#   old 1 "Azure Batch and High Performance Computing" (title)   -> new 1
#   old 2 (blank title placeholder)                              -> REMOVED
#   old 3 "Azure Batch and HPC" / Key learning objectives        -> new 2
#   old 4 "High Performance Computing" / embarrassingly parallel -> new 3
#   old 5 "Azure Batch" / concentrate on the problem...          -> new 4
#   old 6 "Azure High Performance Computing" / MS-MPI...         -> new 5
#   old 7 "Deployment Templates" / easy way to deploy...         -> new 6
#   old 8 "Hands-On Lab"                                         -> new 7
#   old 9 (blank, no shapes)                                     -> new 8

$p = $ppt.ActivePresentation

# Remove the blank slide at position 2 (empty title placeholder only).
$blank = $p.Slides.Item(2)
$blank.Delete()
